$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 6 / Shape 1 (Title placeholder): "User Stories 1" -> "User Stories #1"
# Simple in-place text edit, single run, no formatting change.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange
$sub = $tr.Characters(5, 10)
$sub.Text = " Stories #1"

# ---------------------------------------------------------------------------
# Slide 6 / Shape 3 (textbox): "User Stories 2" -> "User " + "Stories #2"
# The trailing " Stories 2" run is split: a plain space stays in the
# original run, and a brand-new run holding "Stories #2" (sz=3600, bold)
# is created right after it.
# ---------------------------------------------------------------------------
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange
$sub = $tr.Characters(5, 10)
$sub.Text = " Stories #2"
$newPart = $tr.Characters(6, 10)
$newPart.Font.Size = 36
$newPart.Font.Bold = $true

# ---------------------------------------------------------------------------
# Slides 7 & 8: "User Stories N" -> "User Stories " + "#" (bigger) + "N"
# (applies to N = 3, 4, 5, 6)
# ---------------------------------------------------------------------------
function Set-StoriesHash($Shape, $Digit) {
    $tr = $Shape.TextFrame.TextRange
    $sub = $tr.Characters(5, 10)
    $sub.Text = " Stories #$Digit"

    $hashPart = $tr.Characters(14, 1)
    $hashPart.Font.Size = 36
    $hashPart.Font.Bold = $true

    $digitPart = $tr.Characters(15, 1)
    $digitPart.Font.Size = 33
    $digitPart.Font.Bold = $true
}

$s7 = $p.Slides.Item(7)
Set-StoriesHash $s7.Shapes.Item(1) "3"
Set-StoriesHash $s7.Shapes.Item(3) "4"

$s8 = $p.Slides.Item(8)
Set-StoriesHash $s8.Shapes.Item(1) "5"
Set-StoriesHash $s8.Shapes.Item(3) "6"
